$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 56, shifting rows 56:154 down to 57:155
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with data
$ws.Range("A56").Value = 3
$ws.Range("B56").Value = "Femacal de La Calera"
$ws.Range("C56").Value = "Coquimbo"
$ws.Range("D56").Value = 44477
$ws.Range("E56").Value = 5
$ws.Range("F56").Value = 100112001
$ws.Range("G56").Value = "Berenjena"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 130
$ws.Range("K56").Value = 9000
$ws.Range("L56").Value = 9500
$ws.Range("M56").Value = 9269
$ws.Range("N56").Value = '$/caja 60 unidades'
$ws.Range("O56").Value = "Región de Arica y Parinacota"
$ws.Range("P56").Value = 154
$ws.Range("Q56").Value = 60
$ws.Range("R56").Value = "Hortaliza"

# Copy formatting (style) for D56 from D57 (date style) since D column has style s="2"
$ws.Range("D57").Copy()
$ws.Range("D56").PasteSpecial(-4122) # xlPasteFormats
